$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.410281896591187
$ws.Range("B1").Value = 6.431737422943115
$ws.Range("C1").Value = 3.566759824752808
$ws.Range("D1").Value = 1.593715071678162
$ws.Range("E1").Value = 1.122732639312744
